# #3473 replaced two properties that had gaps
# Updates several property rows in the BPS Data sheet: two properties
# (Medstar POB North Tower / 1801 Pennsylvania Ave. area, President Madison
# Apartments, DPW Vehicle Maintenance Facility 2) are replaced by corrected
# records (Medstar POB South Tower, Hampton House, School Without Walls @
# Francis Stevens), and a handful of other cells (owner names, addresses,
# year built, gross area) are corrected for existing rows. The "Year Built"
# column (I) also has its accidental date-style formatting cleared so the
# values display as plain numbers instead of dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Medstar POB North Tower -> Medstar POB South Tower ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Style = "Normal"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319

# --- Row 3: owner name correction ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"
$ws.Range("I3").Style = "Normal"
$ws.Range("I3").Value = 1991

# --- Row 4: address + owner correction ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("I4").Style = "Normal"
$ws.Range("I4").Value = 1991
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# --- Row 5: gross area correction ---
$ws.Range("I5").Style = "Normal"
$ws.Range("I5").Value = 1962
$ws.Range("L5").Value = 58717

# --- Row 6: President Madison Apartments -> Hampton House ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Style = "Normal"
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580

# --- Row 7: postal code + gross area correction ---
$ws.Range("H7").Value = 20005
$ws.Range("I7").Style = "Normal"
$ws.Range("I7").Value = 2004
$ws.Range("L7").Value = 145697

# --- Row 8: address correction ---
$ws.Range("E8").Value = "1428 H ST NW"
$ws.Range("I8").Style = "Normal"
$ws.Range("I8").Value = 1912

# --- Row 9: clear stray date style only ---
$ws.Range("I9").Style = "Normal"
$ws.Range("I9").Value = 1880

# --- Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Style = "Normal"
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991

# Update the sheet's selection to match the saved view state.
$ws.Range("A1:L10").Select()
